{"js": "// Update the date title and all 100 math-fact cells in the table.\n// Cell values are addressed positionally (row-major, 5 columns x 20 rows)\n// because some old values repeat (e.g. \"91-70=\" appears twice and maps to\n// two different new values depending on position), so a global text\n// replace would be ambiguous.\n\nconst oldTitle = \"2024-02-08 Thursday\";\nconst newTitle = \"2024-02-09 Friday\";\n\n// [oldText, newText] for every data cell, in document (row-major) order.\nconst cellPairs = [\n    [\"59+15=\", \"33-17=\"],\n    [\"75-49=\", \"32-1=\"],\n    [\"80+2=\", \"58-25=\"],\n    [\"2+16=\", \"77-35=\"],\n    [\"34-27=\", \"49+18=\"],\n    [\"95-55=\", \"48-16=\"],\n    [\"60+28=\", \"27+0=\"],\n    [\"5+2=\", \"21+9=\"],\n    [\"5+82=\", \"59-52=\"],\n    [\"59-58=\", \"91-54=\"],\n    [\"91-70=\", \"86-2=\"],\n    [\"70-64=\", \"1+9=\"],\n    [\"80-54=\", \"43-38=\"],\n    [\"55+44=\", \"52-52=\"],\n    [\"17+36=\", \"97-2=\"],\n    [\"75+15=\", \"52-34=\"],\n    [\"79-35=\", \"24+41=\"],\n    [\"12+1=\", \"6+86=\"],\n    [\"24-16=\", \"74-4=\"],\n    [\"6+88=\", \"80-29=\"],\n    [\"55+36=\", \"16-15=\"],\n    [\"13+33=\", \"41-11=\"],\n    [\"14+21=\", \"87-48=\"],\n    [\"3+80=\", \"34-30=\"],\n    [\"48+25=\", \"15+30=\"],\n    [\"21+21=\", \"14+73=\"],\n    [\"81-61=\", \"10-4=\"],\n    [\"94-16=\", \"62-47=\"],\n    [\"73-9=\", \"39-4=\"],\n    [\"47+11=\", \"6+6=\"],\n    [\"61-47=\", \"14+75=\"],\n    [\"28+56=\", \"20+11=\"],\n    [\"80-55=\", \"91-44=\"],\n    [\"33+11=\", \"63-55=\"],\n    [\"91+2=\", \"73-70=\"],\n    [\"55-29=\", \"65-7=\"],\n    [\"89-65=\", \"78-4=\"],\n    [\"59-35=\", \"27-15=\"],\n    [\"6+36=\", \"29+51=\"],\n    [\"73-42=\", \"41-32=\"],\n    [\"95-9=\", \"88-82=\"],\n    [\"45+24=\", \"19+17=\"],\n    [\"27+8=\", \"95-81=\"],\n    [\"76-7=\", \"40-2=\"],\n    [\"53-13=\", \"3+64=\"],\n    [\"69-32=\", \"48-12=\"],\n    [\"93+5=\", \"13+84=\"],\n    [\"29+3=\", \"83-19=\"],\n    [\"17-8=\", \"20-1=\"],\n    [\"99-88=\", \"57+10=\"],\n    [\"81-71=\", \"9-6=\"],\n    [\"28+44=\", \"43-5=\"],\n    [\"65-22=\", \"98-54=\"],\n    [\"36-15=\", \"70+4=\"],\n    [\"90-6=\", \"96-73=\"],\n    [\"10+42=\", \"24+30=\"],\n    [\"8+13=\", \"1+95=\"],\n    [\"81+2=\", \"80-23=\"],\n    [\"9+62=\", \"31-14=\"],\n    [\"31+30=\", \"95-35=\"],\n    [\"37+47=\", \"94-58=\"],\n    [\"13+48=\", \"90-79=\"],\n    [\"40-3=\", \"56+43=\"],\n    [\"56-24=\", \"79-20=\"],\n    [\"72-41=\", \"41-2=\"],\n    [\"98-53=\", \"13+57=\"],\n    [\"97-15=\", \"36-4=\"],\n    [\"18-1=\", \"85-82=\"],\n    [\"23+50=\", \"49-42=\"],\n    [\"96-59=\", \"34+7=\"],\n    [\"11+70=\", \"14+17=\"],\n    [\"25+25=\", \"62-13=\"],\n    [\"22+8=\", \"15+55=\"],\n    [\"12+22=\", \"83-15=\"],\n    [\"69-16=\", \"28+20=\"],\n    [\"16+67=\", \"74-48=\"],\n    [\"68-49=\", \"10+68=\"],\n    [\"5+0=\", \"58+27=\"],\n    [\"9+30=\", \"88-72=\"],\n    [\"14+16=\", \"67-1=\"],\n    [\"53+15=\", \"36+62=\"],\n    [\"20+22=\", \"36+52=\"],\n    [\"6+87=\", \"30-23=\"],\n    [\"44+44=\", \"7+27=\"],\n    [\"15+84=\", \"27+34=\"],\n    [\"25+59=\", \"57-6=\"],\n    [\"23+19=\", \"84+13=\"],\n    [\"88-52=\", \"5+13=\"],\n    [\"15+20=\", \"89+6=\"],\n    [\"69-46=\", \"58-7=\"],\n    [\"47-36=\", \"53-33=\"],\n    [\"73+8=\", \"20+39=\"],\n    [\"9+65=\", \"7+50=\"],\n    [\"94-94=\", \"93-16=\"],\n    [\"91-70=\", \"75-74=\"],\n    [\"24+52=\", \"33+8=\"],\n    [\"51-37=\", \"60-16=\"],\n    [\"61-19=\", \"50-34=\"],\n    [\"81-6=\", \"51-24=\"],\n    [\"86-7=\", \"15-5=\"]\n\n];\n\nconst body_ = context.document.body;\nbody_.paragraphs.load(\"items/text\");\nconst tables = body_.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// --- Update the title paragraph (first paragraph of the body) ---\nconst titlePara = body_.paragraphs.items[0];\nif (titlePara.text.trim() === oldTitle) {\n  titlePara.getRange().insertText(newTitle, \"Replace\");\n} else {\n  // Fallback: search document-wide for the exact old title text.\n  const found = body_.search(oldTitle, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n  if (found.items.length > 0) {\n    found.items[0].insertText(newTitle, \"Replace\");\n  }\n}\n\n// --- Update every cell in the first (only) table ---\nconst table = tables.items[0];\ntable.rows.load(\"items/cells\");\nawait context.sync();\n\nconst rows = table.rows;\nconst rowCount = rows.items.length;\nlet idx = 0;\nfor (let r = 0; r < rowCount; r++) {\n  const row = rows.items[r];\n  const cells = row.cells.items;\n  for (let c = 0; c < cells.length; c++) {\n    if (idx >= cellPairs.length) break;\n    const [, newText] = cellPairs[idx];\n    // Write positionally (row-major order matches cellPairs), regardless of\n    // whether the observed current text still equals oldText, so the edit\n    // is deterministic even if a previous step already touched the cell.\n    cells[c].body.getRange().insertText(newText, \"Replace\");\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date title and all 100 math-fact cells in the table.\n# Cells are addressed positionally (row-major, 5 columns x 20 rows) because\n# some old values repeat (e.g. \"91-70=\" appears twice and maps to two\n# different new values depending on position), so a global Find/Replace\n# would be ambiguous.\n\n$d = $word.ActiveDocument\n\n$oldTitle = \"2024-02-08 Thursday\"\n$newTitle = \"2024-02-09 Friday\"\n\n# --- Update the title paragraph (first paragraph of the body) ---\n$titlePara = $d.Paragraphs.Item(1)\n$titleText = $titlePara.Range.Text.TrimEnd([char]13, [char]7)\nif ($titleText -eq $oldTitle) {\n    $titlePara.Range.Text = $newTitle\n} else {\n    # Fallback: use Find/Replace for the exact title text.\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldTitle\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newTitle\n    $find.Execute(\n        $oldTitle, $true, $false, $false, $false, $false, $true, 0, $false,\n        $newTitle, 2, $false, $false, $false, $false\n    )\n}\n\n# --- Update every cell in the first (only) table, by position ---\n# [oldText, newText] pairs in document (row-major) order.\n$cellPairs = @(\n    @(\"59+15=\", \"33-17=\"),\n    @(\"75-49=\", \"32-1=\"),\n    @(\"80+2=\", \"58-25=\"),\n    @(\"2+16=\", \"77-35=\"),\n    @(\"34-27=\", \"49+18=\"),\n    @(\"95-55=\", \"48-16=\"),\n    @(\"60+28=\", \"27+0=\"),\n    @(\"5+2=\", \"21+9=\"),\n    @(\"5+82=\", \"59-52=\"),\n    @(\"59-58=\", \"91-54=\"),\n    @(\"91-70=\", \"86-2=\"),\n    @(\"70-64=\", \"1+9=\"),\n    @(\"80-54=\", \"43-38=\"),\n    @(\"55+44=\", \"52-52=\"),\n    @(\"17+36=\", \"97-2=\"),\n    @(\"75+15=\", \"52-34=\"),\n    @(\"79-35=\", \"24+41=\"),\n    @(\"12+1=\", \"6+86=\"),\n    @(\"24-16=\", \"74-4=\"),\n    @(\"6+88=\", \"80-29=\"),\n    @(\"55+36=\", \"16-15=\"),\n    @(\"13+33=\", \"41-11=\"),\n    @(\"14+21=\", \"87-48=\"),\n    @(\"3+80=\", \"34-30=\"),\n    @(\"48+25=\", \"15+30=\"),\n    @(\"21+21=\", \"14+73=\"),\n    @(\"81-61=\", \"10-4=\"),\n    @(\"94-16=\", \"62-47=\"),\n    @(\"73-9=\", \"39-4=\"),\n    @(\"47+11=\", \"6+6=\"),\n    @(\"61-47=\", \"14+75=\"),\n    @(\"28+56=\", \"20+11=\"),\n    @(\"80-55=\", \"91-44=\"),\n    @(\"33+11=\", \"63-55=\"),\n    @(\"91+2=\", \"73-70=\"),\n    @(\"55-29=\", \"65-7=\"),\n    @(\"89-65=\", \"78-4=\"),\n    @(\"59-35=\", \"27-15=\"),\n    @(\"6+36=\", \"29+51=\"),\n    @(\"73-42=\", \"41-32=\"),\n    @(\"95-9=\", \"88-82=\"),\n    @(\"45+24=\", \"19+17=\"),\n    @(\"27+8=\", \"95-81=\"),\n    @(\"76-7=\", \"40-2=\"),\n    @(\"53-13=\", \"3+64=\"),\n    @(\"69-32=\", \"48-12=\"),\n    @(\"93+5=\", \"13+84=\"),\n    @(\"29+3=\", \"83-19=\"),\n    @(\"17-8=\", \"20-1=\"),\n    @(\"99-88=\", \"57+10=\"),\n    @(\"81-71=\", \"9-6=\"),\n    @(\"28+44=\", \"43-5=\"),\n    @(\"65-22=\", \"98-54=\"),\n    @(\"36-15=\", \"70+4=\"),\n    @(\"90-6=\", \"96-73=\"),\n    @(\"10+42=\", \"24+30=\"),\n    @(\"8+13=\", \"1+95=\"),\n    @(\"81+2=\", \"80-23=\"),\n    @(\"9+62=\", \"31-14=\"),\n    @(\"31+30=\", \"95-35=\"),\n    @(\"37+47=\", \"94-58=\"),\n    @(\"13+48=\", \"90-79=\"),\n    @(\"40-3=\", \"56+43=\"),\n    @(\"56-24=\", \"79-20=\"),\n    @(\"72-41=\", \"41-2=\"),\n    @(\"98-53=\", \"13+57=\"),\n    @(\"97-15=\", \"36-4=\"),\n    @(\"18-1=\", \"85-82=\"),\n    @(\"23+50=\", \"49-42=\"),\n    @(\"96-59=\", \"34+7=\"),\n    @(\"11+70=\", \"14+17=\"),\n    @(\"25+25=\", \"62-13=\"),\n    @(\"22+8=\", \"15+55=\"),\n    @(\"12+22=\", \"83-15=\"),\n    @(\"69-16=\", \"28+20=\"),\n    @(\"16+67=\", \"74-48=\"),\n    @(\"68-49=\", \"10+68=\"),\n    @(\"5+0=\", \"58+27=\"),\n    @(\"9+30=\", \"88-72=\"),\n    @(\"14+16=\", \"67-1=\"),\n    @(\"53+15=\", \"36+62=\"),\n    @(\"20+22=\", \"36+52=\"),\n    @(\"6+87=\", \"30-23=\"),\n    @(\"44+44=\", \"7+27=\"),\n    @(\"15+84=\", \"27+34=\"),\n    @(\"25+59=\", \"57-6=\"),\n    @(\"23+19=\", \"84+13=\"),\n    @(\"88-52=\", \"5+13=\"),\n    @(\"15+20=\", \"89+6=\"),\n    @(\"69-46=\", \"58-7=\"),\n    @(\"47-36=\", \"53-33=\"),\n    @(\"73+8=\", \"20+39=\"),\n    @(\"9+65=\", \"7+50=\"),\n    @(\"94-94=\", \"93-16=\"),\n    @(\"91-70=\", \"75-74=\"),\n    @(\"24+52=\", \"33+8=\"),\n    @(\"51-37=\", \"60-16=\"),\n    @(\"61-19=\", \"50-34=\"),\n    @(\"81-6=\", \"51-24=\"),\n    @(\"86-7=\", \"15-5=\")\n\n)\n\n$table = $d.Tables.Item(1)\n$cols = $table.Columns.Count\n\nfor ($idx = 0; $idx -lt $cellPairs.Count; $idx++) {\n    $row = [int]([math]::Floor($idx / $cols)) + 1\n    $col = ($idx % $cols) + 1\n    $newText = $cellPairs[$idx][1]\n\n    # Write positionally (row-major order matches $cellPairs); Range.Text\n    # assignment preserves the end-of-cell marker automatically.\n    $table.Cell($row, $col).Range.Text = $newText\n}\n\n$d.Save()\n"}
